$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Detaljno budzet": reduce the Researcher effort for the first WP line
# (D14) from 2 to 1 person-month. All dependent formulas (I14, J14, O14,
# Q14, R14 and the row-23 SUM totals) recalculate automatically.
# ---------------------------------------------------------------------------
$wsDetaljno = $wb.Worksheets.Item("Detaljno budzet")
$wsDetaljno.Range("D14").Value = 1

# ---------------------------------------------------------------------------
# Sheet "Travel - budzet": update costs of stay for the existing WP1/WP4/WP5
# trips, then replace the WP7/WP8 trips that used to sit in rows 10-11 with
# a new WP5 (Stuttgart) and WP6 (Lisbon) trip, and move the original WP7
# (Amsterdam) / WP8 (Warsaw) trips down into rows 12-13.
# ---------------------------------------------------------------------------
$wsTravel = $wb.Worksheets.Item("Travel - budzet")

# Row 5 - WP1 / Belgrade -> Athens
$wsTravel.Range("N5").Value = 1260
$wsTravel.Range("O5").Value = 880

# Row 6 - WP4 / Belgrade -> Stuttgart
$wsTravel.Range("N6").Value = 1260
$wsTravel.Range("O6").Value = 880

# Row 7 - WP5 / Belgrade -> Stuttgart
$wsTravel.Range("N7").Value = 1890
$wsTravel.Range("O7").Value = 1320
$wsTravel.Range("P7").Value = 3210

# Row 8 - WP5 / Belgrade -> Stuttgart
$wsTravel.Range("N8").Value = 1890
$wsTravel.Range("O8").Value = 1320
$wsTravel.Range("P8").Value = 3210

# Row 9 - WP5 / Belgrade -> Stuttgart
$wsTravel.Range("N9").Value = 1890
$wsTravel.Range("O9").Value = 1320
$wsTravel.Range("P9").Value = 3210

# Row 10 - becomes a new WP5 / Belgrade -> Stuttgart trip
$wsTravel.Range("B10").Value = "WP5"
$wsTravel.Range("G10").Value = "Stuttgart"
$wsTravel.Range("J10").Value = 3
$wsTravel.Range("K10").ClearContents()
$wsTravel.Range("L10").ClearContents()
$wsTravel.Range("M10").Value = 4
$wsTravel.Range("N10").Value = 1890
$wsTravel.Range("O10").Value = 1320
$wsTravel.Range("P10").Value = 3210

# Row 11 - becomes a new WP6 / Belgrade -> Lisbon trip
$wsTravel.Range("B11").Value = "WP6"
$wsTravel.Range("G11").Value = "Lisbon"
$wsTravel.Range("J11").Value = 2
$wsTravel.Range("K11").ClearContents()
$wsTravel.Range("L11").Value = 3
$wsTravel.Range("M11").Value = 4
$wsTravel.Range("N11").Value = 6250
$wsTravel.Range("O11").Value = 2200
$wsTravel.Range("P11").Value = 8450

# Row 12 - the former WP7 / Belgrade -> Amsterdam trip moves here
$wsTravel.Range("B12").Value = "WP7"
$wsTravel.Range("C12").Value = "ETF"
$wsTravel.Range("D12").Value = "Univerzitet u Beogradu, Elektrotehnički fakultet"
$wsTravel.Range("E12").Value = "SRB"
$wsTravel.Range("F12").Value = "Belgrade"
$wsTravel.Range("G12").Value = "Amsterdam"
$wsTravel.Range("J12").Value = 4
$wsTravel.Range("M12").Value = 5
$wsTravel.Range("N12").Value = 2520
$wsTravel.Range("O12").Value = 2560
$wsTravel.Range("P12").Value = 5080

# Row 13 - the former WP8 / Belgrade -> Warsaw trip moves here
$wsTravel.Range("B13").Value = "WP8"
$wsTravel.Range("C13").Value = "ETF"
$wsTravel.Range("D13").Value = "Univerzitet u Beogradu, Elektrotehnički fakultet"
$wsTravel.Range("E13").Value = "SRB"
$wsTravel.Range("F13").Value = "Belgrade"
$wsTravel.Range("G13").Value = "Warsaw"
$wsTravel.Range("J13").Value = 1
$wsTravel.Range("L13").Value = 1
$wsTravel.Range("M13").Value = 4
$wsTravel.Range("N13").Value = 1260
$wsTravel.Range("O13").Value = 880
$wsTravel.Range("P13").Formula = "=N13+O13"
